$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet
$ws.Name = "RGossF"

# Update a few values in row 13 (precision corrections)
$ws.Range("C13").Value = 0.9862841412949093
$ws.Range("F13").Value = 0.9862841412949093
$ws.Range("L13").Value = 0.990545013099271
$ws.Range("M13").Value = 0.9933429725464397

# Add new row 16 of data
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C16").Value = 1.061031675377856
$ws.Range("D16").Value = 1.233246932143259
$ws.Range("E16").Value = 1.000829744036772
$ws.Range("F16").Value = 1.061031675377856
$ws.Range("G16").Value = 0.8675645298150684
$ws.Range("H16").Value = 1.451569371944928
$ws.Range("I16").Value = 0.962324404813809
$ws.Range("J16").Value = 1.233246932143259
$ws.Range("K16").Value = 1.117038338090016
$ws.Range("L16").Value = 1.089035006733936
$ws.Range("M16").Value = 1.096094443021949

# Match the formatting of column A in row 16 to the rest of column A (bold, bordered, centered style)
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Update the dimension to reflect the new extent
$ws.UsedRange | Out-Null
